$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.325.93"
$ws.Range("E2").Value = "  +4.29%  "
$ws.Range("D3").Value = "2.965.69"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.82%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "2.961.70"
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.14%  "
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").Value = "3.462.45"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").Value = "64.336.69"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.09%  "
$ws.Range("D19").Value = "2.969.74"
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "443.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.675"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.00%  "
$ws.Range("E27").Value = "  +7.65%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.98%  "
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.25%  "
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("E33").Value = "  +2.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.973"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.67%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "48.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "43.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.76%  "
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.292"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "388.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.05%  "
$ws.Range("D46").Value = "2.769.81"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("E47").Value = "  +4.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000222"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.47%  "
$ws.Range("E51").Value = "  +2.48%  "
